$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35, shifting rows 35:177 down to 36:178
$ws.Rows("35:35").Insert()

# Fill in the new row 35 values (same as former row 35 except for the
# values that changed per the new weekly record)
$ws.Range("A35").Value = 8
$ws.Range("B35").Value = "Terminal La Palmera de La Serena"
$ws.Range("C35").Value = "Coquimbo"
$ws.Range("D35").Value = 44659
$ws.Range("D35").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E35").Value = 4
$ws.Range("F35").Value = 100112037
$ws.Range("G35").Value = "Cebollín"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 1600
$ws.Range("K35").Value = 1100
$ws.Range("L35").Value = 1200
$ws.Range("M35").Value = 1150
$ws.Range("N35").Value = "`$/paquete 6 unidades"
$ws.Range("O35").Value = "Provincia del Elquí"
$ws.Range("P35").Value = 192
$ws.Range("Q35").Value = 6
$ws.Range("R35").Value = "Hortaliza"
